# Added 4wk low sales check - refreshed forecast numbers on the
# "Forecast Comparison" sheet (MyForecast / Inventory Coverage /
# Seasonality Index) and the roll-up figures on the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H), Seasonality Index (L) ---

# W10 (row 2)
$ws1.Range("D2").Value = 111
$ws1.Range("H2").Value = 7.86
$ws1.Range("L2").Value = 0.89

# W11 (row 3)
$ws1.Range("D3").Value = 112
$ws1.Range("H3").Value = 6.79
$ws1.Range("L3").Value = 0.99

# W12 (row 4)
$ws1.Range("D4").Value = 113
$ws1.Range("H4").Value = 5.74
$ws1.Range("L4").Value = 1.13

# W13 (row 5)
$ws1.Range("D5").Value = 112
$ws1.Range("H5").Value = 4.79
$ws1.Range("L5").Value = 1.11

# W14 (row 6)
$ws1.Range("D6").Value = 112
$ws1.Range("H6").Value = 3.79
$ws1.Range("L6").Value = 1.13

# W15 (row 7) - forecast unchanged
$ws1.Range("H7").Value = 2.81
$ws1.Range("L7").Value = 1.04

# W16 (row 8)
$ws1.Range("D8").Value = 112
$ws1.Range("H8").Value = 1.79
$ws1.Range("L8").Value = 1.1

# W17 (row 9)
$ws1.Range("D9").Value = 112
$ws1.Range("H9").Value = 0.79
$ws1.Range("L9").Value = 0.88

# W18 (row 10) - inventory coverage unchanged (0)
$ws1.Range("D10").Value = 112
$ws1.Range("L10").Value = 0.9399999999999999

# W19 (row 11) - inventory coverage unchanged (0)
$ws1.Range("D11").Value = 112
$ws1.Range("L11").Value = 1.12

# W20 (row 12) - inventory coverage unchanged (0)
$ws1.Range("D12").Value = 112
$ws1.Range("L12").Value = 0.82

# W21 (row 13) - inventory coverage unchanged (0)
$ws1.Range("D13").Value = 113
$ws1.Range("L13").Value = 1.16

# W22 (row 14) - inventory coverage unchanged (0)
$ws1.Range("D14").Value = 113
$ws1.Range("L14").Value = 1.13

# W23 (row 15) - inventory coverage unchanged (0)
$ws1.Range("D15").Value = 112
$ws1.Range("L15").Value = 0.95

# W24 (row 16) - inventory coverage unchanged (0)
$ws1.Range("D16").Value = 113
$ws1.Range("L16").Value = 1.09

# W25 (row 17) - inventory coverage unchanged (0)
$ws1.Range("D17").Value = 111
$ws1.Range("L17").Value = 0.95

# --- Summary roll-up values (stored as text in the workbook, so keep them text) ---

$ws2.Range("B9").Value  = "'1793"   # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "'895"    # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "'448"    # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "'113"    # Max Forecast
$ws2.Range("B14").Value = "'111"    # Min Forecast
